# modified born position of city
# RelivePos (column E) of the first scene row (villageScene, row 2)
# changes from "0,0,0" to "20,0,-137".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = "20,0,-137"
